$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.817.82"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "1.597.79"
$ws.Range("E3").Value = "  -2.16%  "
$ws.Range("E4").Value = "  +0.12%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "208.51"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.54%  "
$ws.Range("E6").Value = "  +0.11%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.478"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -5.16%  "
$ws.Range("E8").Value = "  -2.55%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.0609"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "17.83"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -3.45%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0787"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("D12").Value = "1.820.61"
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("D13").Value = "1.593.53"
$ws.Range("E13").Value = "  -2.35%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "4.04"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -3.51%  "
$ws.Range("E15").Value = "  -4.42%  "
$ws.Range("D16").Value = "25.823.83"
$ws.Range("E16").Value = "  -0.60%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "60.37"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("D18").Value = "0.0₃0713"
$ws.Range("E18").Value = "  -4.38%  "
$ws.Range("E19").Value = "  +0.02%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "189.10"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.67%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "4.17"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("E22").Value = "  -2.81%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.93"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -3.24%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.128"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.74%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "141.32"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.26%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.71"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.85%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "6.51"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -4.19%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "14.88"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.21%  "
$ws.Range("E30").Value = "  -2.67%  "
$ws.Range("E31").Value = "  -4.84%  "
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("E33").Value = "  -5.11%  "
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("D36").Value = "1.099.24"
$ws.Range("E36").Value = "  -3.11%  "
$ws.Range("E37").Value = "  -3.07%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.797"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -7.90%  "
$ws.Range("E39").Value = "  -2.80%  "
$ws.Range("E40").Value = "  -5.70%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "95.56"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -3.09%  "
$ws.Range("D42").Value = "1.733.34"
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("E43").Value = "  -2.63%  "
$ws.Range("E44").Value = "  -5.10%  "
$ws.Range("E45").Value = "  -11.69%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "52.95"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.81%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.42"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -3.75%  "
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("E51").Value = "  -2.72%  "
